$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 128
$ws.Range("F4").Value = 936
$ws.Range("F5").Value = 976
$ws.Range("F6").Value = 1744
$ws.Range("F7").Value = 386
$ws.Range("F8").Value = 1163
$ws.Range("F13").Value = 49
$ws.Range("F15").Value = 649
$ws.Range("F17").Value = 92
$ws.Range("F21").Value = 111
$ws.Range("F22").Value = 648
$ws.Range("F26").Value = 33
$ws.Range("F29").Value = 130
$ws.Range("F31").Value = 254

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 615
$ws.Range("F11").Value = 116

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 301

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 301
$ws.Range("F3").Value = 128
$ws.Range("F5").Value = 936
$ws.Range("F6").Value = 976
$ws.Range("F7").Value = 1744
$ws.Range("F8").Value = 386
$ws.Range("F9").Value = 1163
$ws.Range("F15").Value = 49
$ws.Range("F17").Value = 649
$ws.Range("F19").Value = 92
$ws.Range("F29").Value = 111
$ws.Range("F30").Value = 648
$ws.Range("F34").Value = 33
$ws.Range("F39").Value = 130
$ws.Range("F41").Value = 254
$ws.Range("F42").Value = 615
$ws.Range("F43").Value = 116
$ws.Range("F44").Value = 116
